# Completed Associated Relation with User Product Cart
# Rewrites the product rows (IDs, names, prices, image URLs, category ids)
# and appends two new rows so the sheet grows from A1:E13 to A1:E15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 17
$ws.Range("B2").Value = "Xiaoni"
$ws.Range("C2").Value = 1200
$ws.Range("D2").Value = "5000/uploads/Hating Game.png"
$ws.Range("E2").Value = 1

# Row 3
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Nokia"
$ws.Range("C3").Value = 1200
$ws.Range("D3").Value = "5000/uploads/Hating Game.png"
$ws.Range("E3").Value = 1

# Row 4
$ws.Range("A4").Value = 16
$ws.Range("B4").Value = "New Data"
$ws.Range("C4").Value = 4000
$ws.Range("D4").Value = "http://localhost:5000/uploads/Emotion Machine.png"
$ws.Range("E4").Value = 5

# Row 5
$ws.Range("A5").Value = 15
$ws.Range("B5").Value = "Xiaoni"
$ws.Range("C5").Value = 1200
$ws.Range("D5").Value = "5000/uploads/Hating Game.png"
$ws.Range("E5").Value = 5

# Row 6
$ws.Range("A6").Value = 14
$ws.Range("B6").Value = "Xiaoni"
$ws.Range("C6").Value = 1200
$ws.Range("D6").Value = "5000/uploads/Hating Game.png"
$ws.Range("E6").Value = 5

# Row 7
$ws.Range("A7").Value = 13
$ws.Range("B7").Value = "Xiaoni"
$ws.Range("C7").Value = 1200
$ws.Range("D7").Value = "5000/uploads/Hating Game.png"
$ws.Range("E7").Value = 5

# Row 8
$ws.Range("A8").Value = 12
$ws.Range("B8").Value = "Xiaoni"
$ws.Range("C8").Value = 1200
$ws.Range("D8").Value = "5000/uploads/Hating Game.png"
$ws.Range("E8").Value = 5

# Row 9
$ws.Range("A9").Value = 11
$ws.Range("B9").Value = "Xiaoni"
$ws.Range("C9").Value = 1200
$ws.Range("D9").Value = "5000/uploads/Hating Game.png"
$ws.Range("E9").Value = 5

# Row 10
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Xiaoni"
$ws.Range("C10").Value = 1200
$ws.Range("D10").Value = "5000/uploads/Hating Game.png"
$ws.Range("E10").Value = 5

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Samsung"
$ws.Range("C11").Value = 1200
$ws.Range("D11").Value = "5000/uploads/Hating Game.png"
$ws.Range("E11").Value = 5

# Row 12 - Price cell is cleared entirely (no value)
$ws.Range("A12").Value = 23
$ws.Range("B12").Value = "Fault"
$ws.Range("C12").ClearContents()
$ws.Range("D12").Value = "http://localhost:5000/uploads/Fault.png"
$ws.Range("E12").Value = 7

# Row 13
$ws.Range("A13").Value = 22
$ws.Range("B13").Value = "Fault"
$ws.Range("C13").Value = 1700
$ws.Range("D13").Value = "http://localhost:5000/uploads/Fault.png"
$ws.Range("E13").Value = 7

# Row 14 (new row)
$ws.Range("A14").Value = 21
$ws.Range("B14").Value = "Fault"
$ws.Range("C14").Value = 1700
$ws.Range("D14").Value = "http://localhost:5000/uploads/Fault.png"
$ws.Range("E14").Value = 7

# Row 15 (new row)
$ws.Range("A15").Value = 20
$ws.Range("B15").Value = "Huwaei"
$ws.Range("C15").Value = 1700
$ws.Range("D15").Value = "http://localhost:5000//uploads/FasterFene.png"
$ws.Range("E15").Value = 7
